$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell 'D2' '26.949.11'
$ws.Range('E2').Value = '  -0.63%  '
Set-TextCell 'D3' '1.562.51'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextCell 'D5' '207.56'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  -0.68%  '
Set-TextCell 'D10' '0.0600'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  -0.43%  '
Set-TextCell 'D12' '1.784.67'
$ws.Range('E12').Value = '  -0.41%  '
Set-TextCell 'D13' '1.564.82'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('E15').Value = '  -0.61%  '
Set-TextCell 'D16' '62.06'
$ws.Range('E16').Value = '  -0.05%  '
Set-TextCell 'D17' '26.957.74'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('E18').Value = '  +0.93%  '
Set-TextCell 'D19' '216.61'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('E22').Value = '  +0.43%  '
Set-TextCell 'D23' '9.21'
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('E24').Value = '  -1.41%  '
Set-TextCell 'D25' '152.54'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  -0.43%  '
Set-TextCell 'D27' '15.08'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  +0.86%  '
Set-TextCell 'D32' '3.23'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('E33').Value = '  +1.45%  '
Set-TextCell 'D34' '1.421.16'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D35' '1.61'
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D36' '1.06'
$ws.Range('E36').Value = '  +9.95%  '
Set-TextCell 'D37' '2.33'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('E41').Value = '  -1.02%  '
Set-TextCell 'D42' '1.00'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('E44').Value = '  +1.80%  '
Set-TextCell 'D45' '64.82'
$ws.Range('E45').Value = '  +0.18%  '
Set-TextCell 'D46' '1.75'
$ws.Range('E46').Value = '  -1.53%  '
Set-TextCell 'D47' '1.698.32'
$ws.Range('E47').Value = '  -0.48%  '
Set-TextCell 'D48' '87.33'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D49' '0.0520'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D50' '0.0960'
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextCell 'D51' '1.00'
$ws.Range('E51').Value = '  -0.05%  '
